$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.730.35"
$ws.Range("D2").ClearFormats()
$ws.Range("D3").Value = "'2.474.45"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  +1.41%  "
$ws.Range("E3").ClearFormats()
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'  -0.20%  "
$ws.Range("E4").ClearFormats()
$ws.Range("E5").Value = "'  +1.55%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'148.43"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  +1.81%  "
$ws.Range("E6").ClearFormats()
$ws.Range("E7").Value = "'  +0.01%  "
$ws.Range("E7").ClearFormats()
$ws.Range("D8").Value = "'0.542"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'  +1.05%  "
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'2.468.88"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'  +0.98%  "
$ws.Range("E9").ClearFormats()
$ws.Range("D10").Value = "'0.113"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'  +0.46%  "
$ws.Range("E10").ClearFormats()
$ws.Range("E11").Value = "'  +0.99%  "
$ws.Range("E11").ClearFormats()
$ws.Range("D12").Value = "'5.28"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'  +0.31%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'0.359"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  +1.53%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'27.25"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  +1.27%  "
$ws.Range("E14").ClearFormats()
$ws.Range("E15").Value = "'  -2.65%  "
$ws.Range("E15").ClearFormats()
$ws.Range("D16").Value = "'2.924.09"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'  +1.60%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'63.789.02"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  +1.95%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'2.473.15"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  +1.41%  "
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'11.54"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'  +2.01%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'7.40"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  +6.33%  "
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'330.97"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  +2.03%  "
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'4.22"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  +0.94%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'2.07"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  +16.77%  "
$ws.Range("E23").ClearFormats()
$ws.Range("D24").Value = "'1.00"
$ws.Range("D24").ClearFormats()
$ws.Range("D25").Value = "'65.81"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "'  -2.20%  "
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'626.27"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  +11.36%  "
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'0.0000105"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  +2.13%  "
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = "'8.68"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  -1.23%  "
$ws.Range("E28").ClearFormats()
$ws.Range("B29").Value = "'Fetch.AI"
$ws.Range("B29").ClearFormats()
$ws.Range("C29").Value = "'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("C29").ClearFormats()
$ws.Range("D29").Value = "'1.54"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  +4.58%  "
$ws.Range("E29").ClearFormats()
$ws.Range("B30").Value = "'WrappedeETH"
$ws.Range("B30").ClearFormats()
$ws.Range("C30").Value = "'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("C30").ClearFormats()
$ws.Range("D30").Value = "'2.595.66"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "'  +1.32%  "
$ws.Range("E30").ClearFormats()
$ws.Range("E31").Value = "'  +0.43%  "
$ws.Range("E31").ClearFormats()
$ws.Range("D32").Value = "'8.39"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "'  -1.00%  "
$ws.Range("E32").ClearFormats()
$ws.Range("D33").Value = "'0.144"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  -3.15%  "
$ws.Range("E33").ClearFormats()
$ws.Range("D34").Value = "'1.92"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  +1.14%  "
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = "'5.26"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  +6.89%  "
$ws.Range("E35").ClearFormats()
$ws.Range("E36").Value = "'  -1.72%  "
$ws.Range("E36").ClearFormats()
$ws.Range("D37").Value = "'0.999"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  -0.01%  "
$ws.Range("E37").ClearFormats()
$ws.Range("E38").Value = "'  +0.04%  "
$ws.Range("E38").ClearFormats()
$ws.Range("E39").Value = "'  +0.02%  "
$ws.Range("E39").ClearFormats()
$ws.Range("D40").Value = "'18.82"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "'  -0.03%  "
$ws.Range("E40").ClearFormats()
$ws.Range("D41").Value = "'148.44"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  -0.93%  "
$ws.Range("E41").ClearFormats()
$ws.Range("B42").Value = "'dogwifhat"
$ws.Range("B42").ClearFormats()
$ws.Range("C42").Value = "'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("C42").ClearFormats()
$ws.Range("D42").Value = "'2.74"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  +12.97%  "
$ws.Range("E42").ClearFormats()
$ws.Range("B43").Value = "'Stacks"
$ws.Range("B43").ClearFormats()
$ws.Range("C43").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("C43").ClearFormats()
$ws.Range("D43").Value = "'2.74"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "'  -1.21%  "
$ws.Range("E43").ClearFormats()
$ws.Range("E44").Value = "'  -0.22%  "
$ws.Range("E44").ClearFormats()
$ws.Range("D45").Value = "'150.65"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "'  +0.93%  "
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'3.79"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  +2.74%  "
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'21.39"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  +3.62%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'0.0544"
$ws.Range("D48").ClearFormats()
$ws.Range("D49").Value = "'0.605"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  +0.41%  "
$ws.Range("E49").ClearFormats()
$ws.Range("D50").Value = "'0.0237"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "'  +2.04%  "
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'0.0920"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  -1.05%  "
$ws.Range("E51").ClearFormats()
